$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styling on column D while forcing text storage
# for numeric-looking price strings (so COM does not coerce them to numbers).
$dStyle = $ws.Range("D2:D51").Style
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '51.101.65'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '3.058.41'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '389.87'
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("D6").Value = '101.59'
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("E7").Value = '  -2.23%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.578'
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").Value = '36.64'
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '0.0846'
$ws.Range("E12").Value = '  -1.73%  '
$ws.Range("D13").Value = '3.536.87'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").Value = '18.23'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").Value = '7.64'
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").Value = '3.053.89'
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("D17").Value = '0.986'
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '10.61'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Value = '51.081.54'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("D21").Value = '0.0₃0953'
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("D22").Value = '12.19'
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").Value = '69.57'
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").Value = '263.52'
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("D25").Value = '3.11'
$ws.Range("E25").Value = '  -1.55%  '
$ws.Range("D26").Value = '7.85'
$ws.Range("E26").Value = '  -6.95%  '
$ws.Range("D27").Value = '26.64'
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("D29").Value = '7.11'
$ws.Range("E29").Value = '  -5.38%  '
$ws.Range("E30").Value = '  -5.81%  '
$ws.Range("D32").Value = '10.43'
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").Value = '0.0483'
$ws.Range("E33").Value = '  +7.85%  '
$ws.Range("D34").Value = '35.47'
$ws.Range("E34").Value = '  +3.95%  '
$ws.Range("D35").Value = '2.07'
$ws.Range("E35").Value = '  -0.64%  '
$ws.Range("D36").Value = '49.98'
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '3.33'
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("E39").Value = '  -2.50%  '
$ws.Range("D40").Value = '129.62'
$ws.Range("E40").Value = '  +1.68%  '
$ws.Range("D41").Value = '16.45'
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").Value = '1.82'
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("D46").Value = '21.61'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("E47").Value = '  +3.91%  '
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").Value = '2.059.44'
$ws.Range("E49").Value = '  +1.80%  '
$ws.Range("D50").Value = '0.0321'
$ws.Range("E50").Value = '  +3.02%  '
$ws.Range("D51").Value = '0.892'
$ws.Range("E51").Value = '  +12.30%  '

# Restore original style/number-format on column D
$ws.Range("D2:D51").Style = $dStyle
